$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B held the scraped data (header value in B1, URLs in B2:B11) while
# column A only held a running index (and had no real content worth keeping).
# The sheet is being collapsed down to a single column: shift everything from
# column B one column to the left, into column A, then drop column B.

# Row 1: move the header value (0) from B1 into A1, preserving B1's formatting.
$ws.Range("B1").Copy($ws.Range("A1"))

# Rows 2-11: move each URL from column B into column A, preserving formatting
# (these cells were unstyled, so A2:A11 end up unstyled too).
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 2).Copy($ws.Cells.Item($r, 1))
}

# Column B is no longer needed now that its contents live in column A.
$ws.Range("B1:B11").Delete()
